# Setting tab TC-01,02,09 added (27/10/25)
# Adds a new worksheet "tc010" after the existing "tc009" sheet and
# populates it with a small field/value table, mirroring the structure
# of the other "Setting" test-case sheets already in the workbook.

$wb = $excel.ActiveWorkbook

# Insert the new sheet immediately after the last existing sheet (tc009)
# so it becomes the new last tab in the workbook.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "tc010"

# Populate the data - set values in an order that introduces the new
# shared strings "Desciption" then "def_value".
$newSheet.Range("A1").Value = "fieldname"
$newSheet.Range("A2").Value = "Desciption"
$newSheet.Range("B1").Value = "def_value"
$newSheet.Range("B2").Value = "Testing"

# Match the final selection/active cell on the new tab.
$newSheet.Range("B2").Select()
